$wb = $excel.ActiveWorkbook

# --- README sheet: drop the vestigial "applyFont" style (s=2) that Excel
#     had been carrying on column A / a handful of B+C cells. Clearing
#     the format collapses those cells back to the workbook default,
#     which is how they ended up with no explicit style in the saved file.
$ws1 = $wb.Worksheets.Item("README")
$readmeCells = $ws1.Range("A2,B2,A3,A4,B4,A5,B5,A6,B6,A7,A8,A9,A10,A11,B11,A12,C12,A13,A14,A15,A16,A17,A18")
foreach ($area in $readmeCells.Areas) {
    $area.ClearFormats()
}

# --- Template sheet: same cleanup for row 2 (the example/data row), plus
#     a new "Term Type" column with a "germplasm passport" example value.
$ws2 = $wb.Worksheets.Item("Template")
$templateRow2 = $ws2.Range("A2:Q2")
foreach ($area in $templateRow2.Areas) {
    $area.ClearFormats()
}

$ws2.Range("T1").Value = "Term Type"
$ws2.Range("T1").Font.Bold = $true

$ws2.Range("T2").Value = "germplasm passport"

# Move the selection/view to reflect where editing happened.
$ws2.Activate()
$ws2.Range("S8").Select()
